$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 5) to the log sheet at 01‏/05‏/2025 02:12:45 م
$ws.Range("A5").Value = "01‏/05‏/2025 02:12:45 م"
$ws.Range("B5").Value = "WCK"
$ws.Range("C5").Value = "C3"
$ws.Range("D5").Value = "الرحلة 2"
$ws.Range("E5").Value = "ايتا"
$ws.Range("F5").Value = "احمد"

# G5 ("23223") looks numeric - prefix with an apostrophe so Excel keeps it
# as text (matching the rest of the column, which is stored as text), then
# drop back to the Normal style so no stray quote-prefix formatting sticks.
$ws.Range("G5").Value = "'23223"
$ws.Range("G5").Style = "Normal"

# H5 is blank in the source row - write it as an explicit empty text value
# (apostrophe-only entry), same text-typed-but-empty shape as H2/H3, then
# clear the resulting quote-prefix style.
$ws.Range("H5").Value = "'"
$ws.Range("H5").Style = "Normal"
